# Remove the "Ver no Jupiter ..." / copyright footer block that follows the
# "Requisitos" section, while keeping the trailing blank paragraph and the
# page-break paragraph that come after it.
#
# Before:
#   ... LOM3018: Introdução à Engenharia de Materiais (Requisito fraco)
#   <blank paragraph>
#   Ver no Jupiter Salvar em pdf Salvar em docx
#   © 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution
#   <blank paragraph>
#   <page-break paragraph>
#
# After:
#   ... LOM3018: Introdução à Engenharia de Materiais (Requisito fraco)
#   <blank paragraph>
#   <page-break paragraph>

$d = $word.ActiveDocument

$startPar = $null
$endPar = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "Ver no Jupiter") {
        $startPar = $p
    }
    if ($t -match "Powered by Jekyll") {
        $endPar = $p
    }
}

if ($startPar -ne $null -and $endPar -ne $null) {
    # Extend the start backwards to also swallow the blank paragraph that
    # immediately precedes "Ver no Jupiter ...".
    $delStart = $startPar.Previous().Range.Start
    $delEnd = $endPar.Range.End

    $range = $d.Range($delStart, $delEnd)
    $range.Delete()
}
